$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A5, B5, C5 text content
$ws.Range("A5").Value = "Кыргыз Республикасы"
$ws.Range("B5").Value = "Кыргызская Республика"
$ws.Range("C5").Value = "Kyrgyz Republic "

# Add new column R: header 2023 and value 53.5
$ws.Range("R4").Value = 2023
$ws.Range("R5").Value = 53.5

# Copy styles from Q4/Q5 to R4/R5 to match existing pattern
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial(-4122) # xlPasteFormats

# Set column widths A:C to 35.85546875
$ws.Range("A1:C1").ColumnWidth = 35.85546875

# Adjust row 5 height
$ws.Range("A5").RowHeight = 21

# Remove selection override - set the selection back to default? (probably not easy via API)
